$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'57.320.93"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -4.96%  "
$ws.Range("D3").Value = "'3.114.03"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -5.96%  "
$ws.Range("D4").Value = "'0.999"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.05%  "
$ws.Range("D5").Value = "'520.03"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -6.81%  "
$ws.Range("D6").Value = "'134.44"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -5.47%  "
$ws.Range("E7").Value = "  -0.13%  "
$ws.Range("D8").Value = "'3.111.99"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -6.10%  "
$ws.Range("D9").Value = "'0.445"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -6.39%  "
$ws.Range("D10").Value = "'7.15"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -9.25%  "
$ws.Range("D11").Value = "'0.108"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -9.25%  "
$ws.Range("D12").Value = "'0.379"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -7.19%  "
$ws.Range("D13").Value = "'3.639.95"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -6.36%  "
$ws.Range("E14").Value = "  -2.51%  "
$ws.Range("D15").Value = "'25.27"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -6.38%  "
$ws.Range("D16").Value = "'3.105.25"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -6.23%  "
$ws.Range("D17").Value = "'57.299.20"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -5.04%  "
$ws.Range("D18").Value = "'0.0000149"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -10.25%  "
$ws.Range("D19").Value = "'5.73"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -7.43%  "
$ws.Range("D20").Value = "'12.88"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -11.00%  "
$ws.Range("D21").Value = "'7.95"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -8.32%  "
$ws.Range("D22").Value = "'342.01"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -8.92%  "
$ws.Range("E23").Value = "  +0.00%  "
$ws.Range("D24").Value = "'68.34"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -7.81%  "
$ws.Range("D25").Value = "'0.502"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -7.75%  "
$ws.Range("D26").Value = "'3.237.10"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -5.92%  "
$ws.Range("B27").Value = "Kaspa"
$ws.Range("C27").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D27").Value = "'0.166"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -4.25%  "
$ws.Range("B28").Value = "Binance-PegBSC-USD"
$ws.Range("C28").Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$ws.Range("D28").Value = "'1.00"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +0.05%  "
$ws.Range("D29").Value = "'0.0₃0927"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -10.58%  "
$ws.Range("E30").Value = "  -0.11%  "
$ws.Range("D31").Value = "'6.73"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -7.75%  "
$ws.Range("D32").Value = "'6.96"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -9.72%  "
$ws.Range("E33").Value = "  -9.38%  "
$ws.Range("D34").Value = "'21.41"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -5.18%  "
$ws.Range("D35").Value = "'1.22"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -4.48%  "
$ws.Range("D36").Value = "'156.59"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -6.02%  "
$ws.Range("D37").Value = "'4.75"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -8.26%  "
$ws.Range("D38").Value = "'6.13"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -9.39%  "
$ws.Range("E39").Value = "  -10.56%  "
$ws.Range("E40").Value = "  -7.12%  "
$ws.Range("D41").Value = "'0.0685"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -7.39%  "
$ws.Range("B42").Value = "OKB"
$ws.Range("C42").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D42").Value = "'40.23"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -4.22%  "
$ws.Range("B43").Value = "RenzoRestakedETH"
$ws.Range("C43").Value = "https://coinranking.com/coin/lKlJ_MC5M+renzorestakedeth-ezeth"
$ws.Range("D43").Value = "'3.128.14"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -6.47%  "
$ws.Range("D44").Value = "'0.679"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -9.77%  "
$ws.Range("E45").Value = "  -8.09%  "
$ws.Range("B46").Value = "ONDO"
$ws.Range("C46").Value = "https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo"
$ws.Range("D46").Value = "'1.05"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -6.37%  "
$ws.Range("B47").Value = "FirstDigitalUSD"
$ws.Range("C47").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D47").Value = "'0.999"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.23%  "
$ws.Range("D48").Value = "'1.45"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -9.37%  "
$ws.Range("D49").Value = "'2.257.44"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -4.61%  "
$ws.Range("D50").Value = "'6.15"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -6.10%  "
$ws.Range("D51").Value = "'19.69"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -7.81%  "
